# The "Info" sheet's row 2 was a placeholder/template row (team name, addr1..addr5,
# discord#1234, "set none if no community") that sat above the real data row
# (Watchers / iaa1... / stars1... / ...). The edit deletes that placeholder row,
# shifting the real data up from row 3 to row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Rows(2).Delete() | Out-Null

# Reflect the resulting selection state (Excel selects the whole row after a
# row deletion that leaves the cursor on the former row 3, now row 2).
$ws.Rows(2).Select() | Out-Null
